$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for columns B, C, D, E so numeric-looking strings
# (e.g. "1.00", "2.049.86") are preserved exactly as text, not coerced to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '37.144.70'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '2.049.86'
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '249.37'
$ws.Range("E5").Value = '  -2.68%  '
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '55.18'
$ws.Range("E8").Value = '  +16.00%  '
$ws.Range("D9").Value = '61.78'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").Value = '0.0757'
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("E12").Value = '  +5.62%  '
$ws.Range("D13").Value = '15.16'
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("D14").Value = '2.344.74'
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").Value = '2.050.45'
$ws.Range("E17").Value = '  -3.44%  '
$ws.Range("D18").Value = '37.070.80'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '72.51'
$ws.Range("E19").Value = '  -2.06%  '
$ws.Range("D20").Value = '0.0₃0888'
$ws.Range("E20").Value = '  +5.18%  '
$ws.Range("D21").Value = '14.37'
$ws.Range("E21").Value = '  +6.58%  '
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("D23").Value = '238.35'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").Value = '170.70'
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").Value = '9.16'
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = '20.35'
$ws.Range("E28").Value = '  -4.71%  '
$ws.Range("D29").Value = '2.02'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("D32").Value = '1.05'
$ws.Range("E32").Value = '  +14.86%  '
$ws.Range("E33").Value = '  +4.28%  '
$ws.Range("D34").Value = '4.39'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("B36").Value = 'Gas'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D36").Value = '19.49'
$ws.Range("E36").Value = '  -21.20%  '
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = '0.0844'
$ws.Range("E38").Value = '  -11.30%  '
$ws.Range("E39").Value = '  -5.22%  '
$ws.Range("D40").Value = '0.111'
$ws.Range("E40").Value = '  +32.30%  '
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").Value = '18.24'
$ws.Range("E42").Value = '  +12.07%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '1.15'
$ws.Range("E44").Value = '  -3.98%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").Value = '4.42'
$ws.Range("E45").Value = '  +64.51%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '97.34'
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").Value = '1.303.29'
$ws.Range("E48").Value = '  -4.47%  '
$ws.Range("D49").Value = '2.37'
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("E51").Value = '  -3.96%  '
